{"js": "// Update the \"worksheet date\" heading and every three-digit \u00f7 one-digit\n// division problem in the table with a freshly generated one (new dividend,\n// divisor, quotient and remainder). Each \"before\" string below is unique in\n// the document, so an exact-text search safely locates the single run that\n// needs to change.\nconst replacements = [\n    [\"2024-11-18 Monday\", \"2024-11-19 Tuesday\"],\n    [\"414\u00f72=207, 0\", \"719\u00f74=179, 3\"],\n    [\"488\u00f74=122, 0\", \"920\u00f72=460, 0\"],\n    [\"519\u00f73=173, 0\", \"112\u00f79=12, 4\"],\n    [\"881\u00f73=293, 2\", \"562\u00f75=112, 2\"],\n    [\"250\u00f76=41, 4\", \"929\u00f73=309, 2\"],\n    [\"612\u00f72=306, 0\", \"468\u00f77=66, 6\"],\n    [\"868\u00f79=96, 4\", \"746\u00f75=149, 1\"],\n    [\"810\u00f77=115, 5\", \"279\u00f75=55, 4\"],\n    [\"134\u00f73=44, 2\", \"522\u00f75=104, 2\"],\n    [\"526\u00f76=87, 4\", \"483\u00f75=96, 3\"],\n    [\"665\u00f73=221, 2\", \"355\u00f75=71, 0\"],\n    [\"642\u00f73=214, 0\", \"546\u00f78=68, 2\"],\n    [\"708\u00f73=236, 0\", \"654\u00f76=109, 0\"],\n    [\"498\u00f76=83, 0\", \"462\u00f79=51, 3\"],\n    [\"822\u00f74=205, 2\", \"935\u00f79=103, 8\"],\n    [\"955\u00f72=477, 1\", \"303\u00f77=43, 2\"],\n    [\"864\u00f78=108, 0\", \"342\u00f79=38, 0\"],\n    [\"129\u00f72=64, 1\", \"135\u00f75=27, 0\"],\n    [\"656\u00f75=131, 1\", \"586\u00f78=73, 2\"],\n    [\"736\u00f74=184, 0\", \"659\u00f75=131, 4\"],\n    [\"782\u00f76=130, 2\", \"205\u00f79=22, 7\"],\n    [\"132\u00f75=26, 2\", \"942\u00f78=117, 6\"],\n    [\"334\u00f79=37, 1\", \"569\u00f72=284, 1\"],\n    [\"337\u00f72=168, 1\", \"830\u00f78=103, 6\"],\n    [\"889\u00f72=444, 1\", \"804\u00f74=201, 0\"]\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n    const results = body.search(before, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    if (results.items.length === 0) {\n        throw new Error(\"Text not found: \" + before);\n    }\n\n    // Replace every hit (expected to be exactly one per \"before\" string).\n    for (const item of results.items) {\n        item.insertText(after, \"Replace\");\n    }\n    await context.sync();\n}\n", "ps1": "# Update the \"worksheet date\" heading and every three-digit \u00f7 one-digit\n# division problem in the table with a freshly generated one (new dividend,\n# divisor, quotient and remainder). Each \"before\" string is unique in the\n# document, so Find/Replace on the whole-document range safely targets the\n# single run that needs to change.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-11-18 Monday\", \"2024-11-19 Tuesday\"),\n    @(\"414\u00f72=207, 0\", \"719\u00f74=179, 3\"),\n    @(\"488\u00f74=122, 0\", \"920\u00f72=460, 0\"),\n    @(\"519\u00f73=173, 0\", \"112\u00f79=12, 4\"),\n    @(\"881\u00f73=293, 2\", \"562\u00f75=112, 2\"),\n    @(\"250\u00f76=41, 4\", \"929\u00f73=309, 2\"),\n    @(\"612\u00f72=306, 0\", \"468\u00f77=66, 6\"),\n    @(\"868\u00f79=96, 4\", \"746\u00f75=149, 1\"),\n    @(\"810\u00f77=115, 5\", \"279\u00f75=55, 4\"),\n    @(\"134\u00f73=44, 2\", \"522\u00f75=104, 2\"),\n    @(\"526\u00f76=87, 4\", \"483\u00f75=96, 3\"),\n    @(\"665\u00f73=221, 2\", \"355\u00f75=71, 0\"),\n    @(\"642\u00f73=214, 0\", \"546\u00f78=68, 2\"),\n    @(\"708\u00f73=236, 0\", \"654\u00f76=109, 0\"),\n    @(\"498\u00f76=83, 0\", \"462\u00f79=51, 3\"),\n    @(\"822\u00f74=205, 2\", \"935\u00f79=103, 8\"),\n    @(\"955\u00f72=477, 1\", \"303\u00f77=43, 2\"),\n    @(\"864\u00f78=108, 0\", \"342\u00f79=38, 0\"),\n    @(\"129\u00f72=64, 1\", \"135\u00f75=27, 0\"),\n    @(\"656\u00f75=131, 1\", \"586\u00f78=73, 2\"),\n    @(\"736\u00f74=184, 0\", \"659\u00f75=131, 4\"),\n    @(\"782\u00f76=130, 2\", \"205\u00f79=22, 7\"),\n    @(\"132\u00f75=26, 2\", \"942\u00f78=117, 6\"),\n    @(\"334\u00f79=37, 1\", \"569\u00f72=284, 1\"),\n    @(\"337\u00f72=168, 1\", \"830\u00f78=103, 6\"),\n    @(\"889\u00f72=444, 1\", \"804\u00f74=201, 0\"),\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #   MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n    #   ReplaceWith, Replace)\n    #   Wrap:    1 = wdFindContinue (search the whole story)\n    #   Replace: 2 = wdReplaceAll\n    $ok = $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $ok) {\n        throw \"Find/Replace failed for: $findText\"\n    }\n}\n"}
